$wb = $excel.ActiveWorkbook

# --- Rename existing sheet to "Display" and add a new "Work" sheet after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Display"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Work"

# ============================================================
# "Work" sheet: raw 2015 CEA data plus the 2017-dollar update
# ============================================================
$ws2.Range("A1").Value = "Cost Type"
$ws2.Range("B1").Value = "Cost"
$ws2.Range("C1").Value = "Estimate for 2015 (CEA)"
$ws2.Range("D1").Value = "Updated Estimates (2017 dollars)"

$ws2.Range("A2").Value = "Mortality"
$ws2.Range("B2").Value = "Opioid-related Death Costs"
$ws2.Range("C2").Value = 431.7
$ws2.Range("D2").Formula = "=C2*B16*1.04"

$ws2.Range("A3").Value = "Non-Mortality"
$ws2.Range("B3").Value = "Health Care Costs"
$ws2.Range("C3").Value = 36.6
$ws2.Range("D3").Formula = "=B14*B8/1000"

$ws2.Range("A4").Value = "Non-Mortality"
$ws2.Range("B4").Value = "Productivity Costs"
$ws2.Range("C4").Value = 25.9
$ws2.Range("D4").Formula = "=B14*B9/1000"

$ws2.Range("A5").Value = "Non-Mortality"
$ws2.Range("B5").Value = "Criminal Justice Costs"
$ws2.Range("C5").Value = 9.7
$ws2.Range("D5").Formula = "=B14*B10/1000"

$ws2.Range("A7").Value = "Weights"

$ws2.Range("A8").Value = "Health Care Costs"
$ws2.Range("B8").Formula = "=C3/(C3+C4+C5)"

$ws2.Range("A9").Value = "Productivity Costs"
$ws2.Range("B9").Formula = "=C4/(C3+C4+C5)"

$ws2.Range("A10").Value = "Criminal Justice Costs"
$ws2.Range("B10").Formula = "=C5/(C3+C4+C5)"

$ws2.Range("A12").Value = "Non-productivity pp (2015 dollars)"
$ws2.Range("B12").Formula = "=(56990*1.01)/1.9"

$ws2.Range("A13").Value = "Non-productivity pp (2017 dollars)"
$ws2.Range("B13").Formula = "=B12*1.04"

$ws2.Range("A14").Value = "Total non-productivity (2017 Millions)"
$ws2.Range("B14").Formula = "=B13*2.1/1000"

$ws2.Range("A16").Value = "Percent Cost Change"
$ws2.Range("B16").Value = 1.558

# ============================================================
# "Display" sheet: headline table shown to readers
# ============================================================
$ws1.Range("C1").Value = "2015 CEA Estimate (2017 Dollars)"
$ws1.Range("D1").Value = "Updated Updated 2017 Estimates "

$ws1.Range("C2").Formula = "=431.7*1.04"
$ws1.Range("D2").Value = 699.49214400000005

$ws1.Range("C3").Value = 36.6
$ws1.Range("D3").Formula = "=33539.9917667299/1000"

$ws1.Range("C4").Value = 25.9
$ws1.Range("D4").Formula = "=23734.5843376586/1000"

$ws1.Range("C5").Value = 9.7
$ws1.Range("D5").Formula = "=8889.01421140108/1000"

# --- Selections matching the saved view state ---
$ws2.Activate()
[void]$ws2.Range("A13").Select()
$ws1.Activate()
[void]$ws1.Range("D6").Select()
